# Update the "Support" sheet's funding table (Table13): tweak the wording of
# several "Program Funding Sources" entries, and add two new money columns
# ("State AI/GGEE Estimate" and "District/City Estimate") with per-district
# dollar estimates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Support")
$tbl = $ws.ListObjects.Item(1)

# --- 1) Add the two new table columns ---------------------------------------
$colD = $tbl.ListColumns.Add()
$ws.Cells.Item(1, 4).Value2 = "State AI/GGEE Estimate"

$colE = $tbl.ListColumns.Add()
$ws.Cells.Item(1, 5).Value2 = "District/City Estimate"

# Match the header formatting (border/shading) used by the rest of row 1.
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Range("D1:E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2) Re-word several "Program Funding Sources" values in column B -------
# (row 3 / Brevard's "UF Donor Funding" is untouched)
$ws.Cells.Item(5, 2).Value2  = "GGEE, State AI"                # Miami Dade County Public Schools
$ws.Cells.Item(4, 2).Value2  = "City, State AI, GGEE"          # City of Rivieria Beach
$ws.Cells.Item(2, 2).Value2  = "Local, State AI, GGEE"         # Alachua County Schools
$ws.Cells.Item(6, 2).Value2  = "District, State AI, GGEE"      # Orange County Public Schools
$ws.Cells.Item(7, 2).Value2  = "District, State AI, GGEE"      # Pinellas County Schools
$ws.Cells.Item(8, 2).Value2  = "District, State AI, GGEE"      # Santa Rosa County District Schools
$ws.Cells.Item(9, 2).Value2  = "District, State AI, GGEE"      # Sarasota County Schools
$ws.Cells.Item(10, 2).Value2 = "District, State AI, GGEE"      # School District of Palm Beach County

# --- 3) "State AI/GGEE Estimate" (column D) values ---------------------------
# Applying the number format first (in this order) + then the value mirrors
# how the workbook's cellXfs end up allocated.
$ws.Range("D2").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("E2").NumberFormat = '"$"#,##0.00'
$ws.Range("D5").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

$ws.Range("D3").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("D4").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("D6:D10").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'
$ws.Range("E3:E10").NumberFormat = '"$"#,##0.00'

$ws.Cells.Item(2, 4).Value2  = 450
$ws.Cells.Item(3, 4).Value2  = 4300
$ws.Cells.Item(4, 4).Value2  = 2100
$ws.Cells.Item(5, 4).Value2  = 10342
$ws.Cells.Item(6, 4).Value2  = 6700
$ws.Cells.Item(7, 4).Value2  = 3700
$ws.Cells.Item(8, 4).Value2  = 23000
$ws.Cells.Item(9, 4).Value2  = 2000
$ws.Cells.Item(10, 4).Value2 = 3300

# --- 4) "District/City Estimate" (column E) values --------------------------
# This column started life as a calculated column (=SUM of two helper
# columns); those helper columns were later removed, leaving a #REF! formula
# behind, and the cells were then overwritten with the actual numbers.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 26).Value2 = 100
    $ws.Cells.Item($r, 27).Value2 = 200
}
$colE.DataBodyRange.Formula = "=SUM(Z2,AA2)"
$ws.Columns.Item(26).Delete() | Out-Null
$ws.Columns.Item(26).Delete() | Out-Null

$ws.Cells.Item(2, 5).Value2  = 4615
$ws.Cells.Item(3, 5).Value2  = 5900
$ws.Cells.Item(4, 5).Value2  = 8000
$ws.Cells.Item(5, 5).Value2  = 0
$ws.Cells.Item(6, 5).Value2  = 12700
$ws.Cells.Item(7, 5).Value2  = 8100
$ws.Cells.Item(8, 5).Value2  = 18400
$ws.Cells.Item(9, 5).Value2  = 4600
$ws.Cells.Item(10, 5).Value2 = 2100

# --- 5) Stray formatted-but-empty cell below the table ----------------------
$ws.Range("B15").NumberFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# --- 6) Column widths (best-fit for the re-worded/new columns) -------------
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# --- 7) Leave the selection where the editing session ended ----------------
$ws.Range("F2").Select() | Out-Null
